$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.09476902780640728
$ws.Range("A2").Value = -0.0059999999801654269
$ws.Range("A3").Value = -0.003999999984117153
$ws.Range("A4").Value = -0.0079999999691295898
$ws.Range("A5").Value = -0.0029999999834817714
$ws.Range("A6").Value = -0.0019999999824094061
$ws.Range("A7").Value = -0.0099999999566602327
$ws.Range("A8").Value = 0.038813731292671516
$ws.Range("A9").Value = -0.0019999999814288572
$ws.Range("A10").Value = -0.001999999979984679
$ws.Range("A11").Value = -0.002999999976594836
$ws.Range("A12").Value = -0.0034999999745388699
$ws.Range("A13").Value = -0.0034999999731288867
$ws.Range("A14").Value = -0.0079999999581437109
$ws.Range("A15").Value = -0.00099999998070199325
$ws.Range("A16").Value = -0.0019999999772082333
$ws.Range("A17").Value = -0.0019999999769462207
$ws.Range("A18").Value = -0.0039999999703317357
$ws.Range("A19").Value = 0.031303192301287641
$ws.Range("A20").Value = -0.0039999999858668644
$ws.Range("A21").Value = -0.0039999999857274204
$ws.Range("A22").Value = -0.0039999999856084045
$ws.Range("A23").Value = -0.0049999999811429774
$ws.Range("A24").Value = 0.039639428746715666
$ws.Range("A25").Value = -0.019999999932340806
$ws.Range("A26").Value = -0.0024999999798627215
$ws.Range("A27").Value = -0.0024999999784420801
$ws.Range("A28").Value = -0.001999999974423794
$ws.Range("A29").Value = -0.0069999999540701907
$ws.Range("A30").Value = -0.059999999779500879
$ws.Range("A31").Value = -0.0069999999506009658
$ws.Range("A32").Value = -0.0099999999406499285
$ws.Range("A33").Value = -0.0039999999601025849
